$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.459772186609131
$ws.Range("D2").Value = 10.21205681273083
$ws.Range("E2").Value = 8.665687138136612
$ws.Range("F2").Value = 76.12425560518976
$ws.Range("G2").Value = 3.823389022674342
$ws.Range("I2").Value = 18.72193701677038
$ws.Range("J2").Value = 8.56099746325058
$ws.Range("M2").Value = 38.84783946225661
$ws.Range("N2").Value = 17.68517234655202
$ws.Range("B3").Value = 6.300252423865869
$ws.Range("D3").Value = 9.842422498270416
$ws.Range("E3").Value = 8.288230836267784
$ws.Range("F3").Value = 75.26105913222864
$ws.Range("G3").Value = 3.835964348786071
$ws.Range("I3").Value = 18.85249620671442
$ws.Range("J3").Value = 8.57712780139958
$ws.Range("M3").Value = 38.14858128977712
$ws.Range("N3").Value = 17.7438870098695
$ws.Range("B4").Value = 6.200847697060209
$ws.Range("D4").Value = 9.614270175950566
$ws.Range("E4").Value = 8.047090501492411
$ws.Range("F4").Value = 74.76323092004469
$ws.Range("G4").Value = 3.844017905661135
$ws.Range("I4").Value = 18.93634076625256
$ws.Range("J4").Value = 8.587637973490757
$ws.Range("M4").Value = 37.72700855324533
$ws.Range("N4").Value = 17.78221176998855
$ws.Range("B5").Value = 6.160029389850452
$ws.Range("D5").Value = 9.52114507839479
$ws.Range("E5").Value = 7.946525759289502
$ws.Range("F5").Value = 74.56853731995481
$ws.Range("G5").Value = 3.847384255952396
$ws.Range("I5").Value = 18.97143980915805
$ws.Range("J5").Value = 8.592073617087271
$ws.Range("M5").Value = 37.55735586341757
$ws.Range("N5").Value = 17.79840483747738
$ws.Range("B6").Value = 6.153234530243961
$ws.Range("D6").Value = 9.505677172497924
$ws.Range("E6").Value = 7.929690242500087
$ws.Range("F6").Value = 74.53670441852758
$ws.Range("G6").Value = 3.847948363379031
$ws.Range("I6").Value = 18.97732443001371
$ws.Range("J6").Value = 8.592819380542863
$ws.Range("M6").Value = 37.52931990626644
$ws.Range("N6").Value = 17.80112854456142
$ws.Range("B7").Value = 6.200298385375991
$ws.Range("D7").Value = 9.613014665024302
$ws.Range("E7").Value = 8.045743467217035
$ws.Range("F7").Value = 74.7605720161886
$ws.Range("G7").Value = 3.844062962264761
$ws.Range("I7").Value = 18.93681034365292
$ws.Range("J7").Value = 8.587697175650074
$ws.Range("M7").Value = 37.72471164192183
$ws.Range("N7").Value = 17.78242782066656
$ws.Range("B8").Value = 6.405105720082028
$ws.Range("D8").Value = 10.0849411664045
$ws.Range("E8").Value = 8.537520395928924
$ws.Range("F8").Value = 75.82000346533452
$ws.Range("G8").Value = 3.82765659160292
$ws.Range("I8").Value = 18.76619364762735
$ws.Range("J8").Value = 8.56643356737453
$ws.Range("M8").Value = 38.60523427583543
$ws.Range("N8").Value = 17.70494747443444
$ws.Range("B9").Value = 6.792818132609564
$ws.Range("D9").Value = 10.99494016168574
$ws.Range("E9").Value = 9.425311383280073
$ws.Range("F9").Value = 78.14826900908032
$ws.Range("G9").Value = 3.798077619559462
$ws.Range("I9").Value = 18.46053347154428
$ws.Range("J9").Value = 8.529533043521287
$ws.Range("M9").Value = 40.38537790089354
$ws.Range("N9").Value = 17.57088266881241
$ws.Range("B10").Value = 7.066292173849639
$ws.Range("D10").Value = 11.64702635081214
$ws.Range("E10").Value = 10.02865798584999
$ws.Range("F10").Value = 80.00448218189523
$ws.Range("G10").Value = 3.777866436114298
$ws.Range("I10").Value = 18.25318623623212
$ws.Range("J10").Value = 8.505330956493545
$ws.Range("M10").Value = 41.71491854631369
$ws.Range("N10").Value = 17.48304770267812
$ws.Range("B11").Value = 7.187689230625222
$ws.Range("D11").Value = 11.93886587347078
$ws.Range("E11").Value = 10.29220518057061
$ws.Range("F11").Value = 80.87875207246536
$ws.Range("G11").Value = 3.768988296729478
$ws.Range("I11").Value = 18.16250729425039
$ws.Range("J11").Value = 8.494949204944694
$ws.Range("M11").Value = 42.32224158556199
$ws.Range("N11").Value = 17.44535293051805
$ws.Range("B12").Value = 7.233186189773997
$ws.Range("D12").Value = 12.04860271410743
$ws.Range("E12").Value = 10.39041620188743
$ws.Range("F12").Value = 81.21394606398596
$ws.Range("G12").Value = 3.765670672154672
$ws.Range("I12").Value = 18.1286861874119
$ws.Range("J12").Value = 8.49110797765165
$ws.Range("M12").Value = 42.55240316113682
$ws.Range("N12").Value = 17.43139986076054
$ws.Range("B13").Value = 7.223409278416884
$ws.Range("D13").Value = 12.02500470019029
$ws.Range("E13").Value = 10.36933558656539
$ws.Range("F13").Value = 81.14157521233146
$ws.Range("G13").Value = 3.766383228315733
$ws.Range("I13").Value = 18.13594728615396
$ws.Range("J13").Value = 8.491931249354492
$ws.Range("M13").Value = 42.50282842554729
$ws.Range("N13").Value = 17.43439068875562
$ws.Range("B14").Value = 7.191442004112103
$ws.Range("D14").Value = 11.94791000192873
$ws.Range("E14").Value = 10.30031698320866
$ws.Range("F14").Value = 80.90624695112037
$ws.Range("G14").Value = 3.768714471743202
$ws.Range("I14").Value = 18.15971449003323
$ws.Range("J14").Value = 8.494631379743351
$ws.Range("M14").Value = 42.34117447137907
$ws.Range("N14").Value = 17.44419859128531
$ws.Range("B15").Value = 7.17179830958832
$ws.Range("D15").Value = 11.900583801349
$ws.Range("E15").Value = 10.25783376615848
$ws.Range("F15").Value = 80.76263398542763
$ws.Range("G15").Value = 3.770148165616478
$ws.Range("I15").Value = 18.17433971118345
$ws.Range("J15").Value = 8.49629701846999
$ws.Range("M15").Value = 42.24217524009579
$ws.Range("N15").Value = 17.45024791028263
$ws.Range("B16").Value = 7.058294485798442
$ws.Range("D16").Value = 11.62784972346312
$ws.Range("E16").Value = 10.01121313100099
$ws.Range("F16").Value = 79.94793336015285
$ws.Range("G16").Value = 3.778452908708527
$ws.Range("I16").Value = 18.25918505691083
$ws.Range("J16").Value = 8.506022048110184
$ws.Range("M16").Value = 41.67526350384242
$ws.Range("N16").Value = 17.48555630789995
$ws.Range("B17").Value = 6.987862957555817
$ws.Range("D17").Value = 11.45924179553406
$ws.Range("E17").Value = 9.857107811236926
$ws.Range("F17").Value = 79.45567304102698
$ws.Range("G17").Value = 3.783627764747643
$ws.Range("I17").Value = 18.31216340453538
$ws.Range("J17").Value = 8.512148730758708
$ws.Range("M17").Value = 41.32799359526228
$ws.Range("N17").Value = 17.5077932666172
$ws.Range("B18").Value = 6.947071733692766
$ws.Range("D18").Value = 11.36181574512659
$ws.Range("E18").Value = 9.767443537778174
$ws.Range("F18").Value = 79.17535983917239
$ws.Range("G18").Value = 3.786634020040647
$ws.Range("I18").Value = 18.34297860345963
$ws.Range("J18").Value = 8.515731747111188
$ws.Range("M18").Value = 41.1284994129725
$ws.Range("N18").Value = 17.5207965204593
$ws.Range("B19").Value = 6.933213592571158
$ws.Range("D19").Value = 11.3287549898254
$ws.Range("E19").Value = 9.736909036996549
$ws.Range("F19").Value = 79.08094031805068
$ws.Range("G19").Value = 3.787657040360081
$ws.Range("I19").Value = 18.35347129958487
$ws.Range("J19").Value = 8.516955052692822
$ws.Range("M19").Value = 41.06100201969339
$ws.Range("N19").Value = 17.52523593021055
$ws.Range("B20").Value = 6.9953898911259
$ws.Range("D20").Value = 11.47723735067529
$ws.Range("E20").Value = 9.873619026450854
$ws.Range("F20").Value = 79.50778413929963
$ws.Range("G20").Value = 3.783073814669226
$ws.Range("I20").Value = 18.30648827294965
$ws.Range("J20").Value = 8.511490418521154
$ws.Range("M20").Value = 41.36493690794151
$ws.Range("N20").Value = 17.50540407439686
$ws.Range("B21").Value = 7.200844726471033
$ws.Range("D21").Value = 11.97057628403428
$ws.Range("E21").Value = 10.32063265225087
$ws.Range("F21").Value = 80.97525792172183
$ws.Range("G21").Value = 3.768028534876174
$ws.Range("I21").Value = 18.1527195106548
$ws.Range("J21").Value = 8.493835841581145
$ws.Range("M21").Value = 42.38865261054476
$ws.Range("N21").Value = 17.44130909077122
$ws.Range("B22").Value = 7.332345037561152
$ws.Range("D22").Value = 12.28844460397234
$ws.Range("E22").Value = 10.60351787702461
$ws.Range("F22").Value = 81.95831180426943
$ws.Range("G22").Value = 3.758453442913724
$ws.Range("I22").Value = 18.05523304693448
$ws.Range("J22").Value = 8.482822752558864
$ws.Range("M22").Value = 43.05869195877683
$ws.Range("N22").Value = 17.40128906955747
$ws.Range("B23").Value = 7.262427511110167
$ws.Range("D23").Value = 12.11923466727079
$ws.Range("E23").Value = 10.45338899465813
$ws.Range("F23").Value = 81.43150018224674
$ws.Range("G23").Value = 3.763540632430908
$ws.Range("I23").Value = 18.10699033197873
$ws.Range("J23").Value = 8.488652644048294
$ws.Range("M23").Value = 42.70104616541318
$ws.Range("N23").Value = 17.42247881523352
$ws.Range("B24").Value = 6.991987893253575
$ws.Range("D24").Value = 11.46910308317748
$ws.Range("E24").Value = 9.866157628271324
$ws.Range("F24").Value = 79.4842163084158
$ws.Range("G24").Value = 3.783324158354523
$ws.Range("I24").Value = 18.30905288668452
$ws.Range("J24").Value = 8.511787852388334
$ws.Range("M24").Value = 41.34823433693503
$ws.Range("N24").Value = 17.50648354597071
$ws.Range("B25").Value = 6.689733737416992
$ws.Range("D25").Value = 10.75119554626008
$ws.Range("E25").Value = 9.193600768993196
$ws.Range("F25").Value = 77.49225786692405
$ws.Range("G25").Value = 3.805807995232185
$ws.Range("I25").Value = 18.54016839381283
$ws.Range("J25").Value = 8.539003699213083
$ws.Range("M25").Value = 39.89921967282727
$ws.Range("N25").Value = 17.60526114694385
